$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1240.2
$ws.Range("I4").Value = 109.333336
$ws.Range("J4").Value = 2936.5
$ws.Range("K4").Value = 109.333336
$ws.Range("L4").Value = 2936.5
$ws.Range("M4").Value = 4.666663999999997
$ws.Range("N4").Value = -3164.5
$ws.Range("H6").Value = 1131
$ws.Range("I6").Value = 535.1667
$ws.Range("J6").Value = 2024.75
$ws.Range("K6").Value = 1605.5001
$ws.Range("L6").Value = 6074.25
$ws.Range("M6").Value = -1493.5001
$ws.Range("N6").Value = -6298.25
$ws.Range("H20").Value = 900
$ws.Range("I20").Value = 900
$ws.Range("K20").Value = 900
$ws.Range("M20").Value = -670
$ws.Range("H35").Value = 900
$ws.Range("I35").Value = 900
$ws.Range("K35").Value = 900
$ws.Range("M35").Value = -521
$ws.Range("H76").Value = 4416.8335
$ws.Range("I76").Value = 3500.25
$ws.Range("K76").Value = 3500.25
$ws.Range("M76").Value = -3185.25
$ws.Range("H79").Value = 4416.8335
$ws.Range("I79").Value = 3500.25
$ws.Range("K79").Value = 3500.25
$ws.Range("M79").Value = -2408.25
$ws.Range("H137").Value = 3535.647
$ws.Range("J137").Value = 3829
$ws.Range("L137").Value = 11487
$ws.Range("N137").Value = -16587
$ws.Range("H141").Value = 3279.8572
$ws.Range("I141").Value = 1410.8
$ws.Range("K141").Value = 4232.4
$ws.Range("M141").Value = 947.6000000000004

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").ClearContents()
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()
$ws.Range("H138").Value = 100000
$ws.Range("J138").Value = 100000
$ws.Range("L138").Value = 100000
$ws.Range("N138").Value = -110280

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 131.23077
$ws.Range("I7").Value = 114.75
$ws.Range("J7").Value = 157.6
$ws.Range("K7").Value = 114.75
$ws.Range("L7").Value = 157.6
$ws.Range("M7").Value = -1.75
$ws.Range("N7").Value = -383.6
$ws.Range("H19").Value = 108.333336
$ws.Range("I19").Value = 108.333336
$ws.Range("K19").Value = 108.333336
$ws.Range("M19").Value = 61.666664
$ws.Range("H24").Value = 108.333336
$ws.Range("I24").Value = 108.333336
$ws.Range("K24").Value = 108.333336
$ws.Range("M24").Value = 61.666664
$ws.Range("H88").Value = 10976.4
$ws.Range("J88").Value = 12085.25
$ws.Range("L88").Value = 12085.25
$ws.Range("N88").Value = -12897.25
$ws.Range("H91").Value = 10976.4
$ws.Range("J91").Value = 12085.25
$ws.Range("L91").Value = 12085.25
$ws.Range("N91").Value = -14893.25
$ws.Range("H96").Value = 9853.416999999999
$ws.Range("J96").Value = 9853.416999999999
$ws.Range("L96").Value = 9853.416999999999
$ws.Range("N96").Value = -15345.417
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 25.7
$ws.Range("I2").Value = 26.11111
$ws.Range("J2").Value = 22
$ws.Range("K2").Value = 156.66666
$ws.Range("L2").Value = 132
$ws.Range("M2").Value = -43.66666000000001
$ws.Range("N2").Value = -358
$ws.Range("H16").Value = 1686.6666
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 2525
$ws.Range("K16").Value = 30
$ws.Range("L16").Value = 7575
$ws.Range("M16").Value = 143
$ws.Range("N16").Value = -7921
$ws.Range("H60").Value = 1058.0952
$ws.Range("I60").Value = 299.6154
$ws.Range("J60").Value = 2290.625
$ws.Range("K60").Value = 898.8462000000001
$ws.Range("L60").Value = 6871.875
$ws.Range("M60").Value = -647.8462000000001
$ws.Range("N60").Value = -7373.875
$ws.Range("H128").Value = 609996.8
$ws.Range("I128").Value = 609996.8
$ws.Range("K128").Value = 1829990.4
$ws.Range("M128").Value = -1825010.4
$ws.Range("H129").Value = 2092.6155
$ws.Range("I129").Value = 651.125
$ws.Range("J129").Value = 4399
$ws.Range("K129").Value = 1953.375
$ws.Range("L129").Value = 13197
$ws.Range("M129").Value = 3046.625
$ws.Range("N129").Value = -23197

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3500
$ws.Range("I5").Value = 3500
$ws.Range("K5").Value = 3500
$ws.Range("M5").Value = -3388
$ws.Range("H9").Value = 206.4
$ws.Range("I9").Value = 233
$ws.Range("J9").Value = 100
$ws.Range("K9").Value = 233
$ws.Range("L9").Value = 100
$ws.Range("M9").Value = -63
$ws.Range("N9").Value = -440

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H30").Value = 432.83334
$ws.Range("I30").Value = 432.83334
$ws.Range("K30").Value = 432.83334
$ws.Range("M30").Value = -324.83334
$ws.Range("H35").Value = 829.4
$ws.Range("I35").Value = 829.4
$ws.Range("J35").Value = 0
$ws.Range("K35").Value = 829.4
$ws.Range("L35").Value = 0
$ws.Range("M35").Value = -493.4
$ws.Range("N35").ClearContents()
$ws.Range("H46").Value = 5408.8823
$ws.Range("I46").Value = 3832.2856
$ws.Range("J46").Value = 6512.5
$ws.Range("K46").Value = 3832.2856
$ws.Range("L46").Value = 6512.5
$ws.Range("M46").Value = -3644.2856
$ws.Range("N46").Value = -6888.5
$ws.Range("H81").Value = 19999
$ws.Range("J81").Value = 19999
$ws.Range("L81").Value = 19999
$ws.Range("N81").Value = -21995
$ws.Range("H84").Value = 19999
$ws.Range("J84").Value = 19999
$ws.Range("L84").Value = 59997
$ws.Range("N84").Value = -69981

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3518.1333
$ws.Range("I132").Value = 1782.3334
$ws.Range("J132").Value = 4675.3335
$ws.Range("K132").Value = 5347.0002
$ws.Range("L132").Value = 14026.0005
$ws.Range("M132").Value = -2817.0002
$ws.Range("N132").Value = -19086.0005
